$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '64.730.37'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +2.03%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.632.54'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +2.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '592.56'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.60%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '154.69'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.84%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.590'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.46%  '
$ws.Range('E9').Value = '  +5.52%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.396'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +3.09%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.76'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.64%  '
$ws.Range('E12').Value = '  +1.69%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '28.83'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +4.63%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000185'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +18.39%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.103.20'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.97%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.694.91'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.28%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.589.45'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +1.28%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '12.51'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +2.67%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.77'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.85%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '349.95'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.95%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.24'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +5.58%  '
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '67.88'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.98%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.69'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.34%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.44'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +3.24%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.62'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.80%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.02'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.64%  '
$ws.Range('B28').Value = 'Aptos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.07'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.163'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0₃0928'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +7.75%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.08'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.10%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '509.64'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -8.12%  '
$ws.Range('E33').Value = '  +0.32%  '
$ws.Range('E34').Value = '  +6.59%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '6.19'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +2.00%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.422'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +2.14%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '164.64'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.23%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '20.05'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +2.55%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.99'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +4.15%  '
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '42.21'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +6.36%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '163.73'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.21%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '4.06'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.87%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0609'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +3.85%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '22.70'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.16%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.17'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +3.26%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.644'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +2.62%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0252'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.64%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0977'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.42%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '19.22'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.48%  '

Write-Output "done"